# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, bordered, centered) from the last
# existing header cell (AC1) onto the three new header cells so they
# reuse the same style as the rest of row 1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (59 wins, 103 losses, 0 ties) for every
# player row.
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 59
    $ws.Cells.Item($r, 31).Value = 103
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-53"
